# Atualização de bases das ligas, do dia: 21-04-2024 às 13:33
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Swap the full content (columns B..AC) of rows 70 and 71 - column A (the
#    sequential id) is left untouched.
# ---------------------------------------------------------------------------
$row70 = @("6779644","FC Nordsjaelland","Odense BK",0,1,"A",1.5,4.2,6,1.333,5,9,-1.5,1.85,2,3.5,2.025,1.825,-1,-1,8,-1,1,-1,0.825)
$row71 = @("6779645","Vejle","Hvidovre IF",3,1,"H",1.833,3.6,4.2,1.8,3.6,4.5,-0.5,1.825,2.025,2.5,2.025,1.825,0.8,-1,-1,0.825,-1,1.025,-1)

# Columns B..AC in order, skipping C/D (Div / Div Original Name, unchanged)
# and E (Date, unchanged). The arrays above hold: B,F,G,H,I,J,K,L,M,N,O,P,Q,R,S,T,U,V,W,X,Y,Z,AA,AB,AC
$cols = @(2,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28,29)

for ($i = 0; $i -lt $cols.Count; $i++) {
    $ws.Cells.Item(70, $cols[$i]).Value = $row71[$i]
    $ws.Cells.Item(71, $cols[$i]).Value = $row70[$i]
}

# ---------------------------------------------------------------------------
# 2) Swap the full content (columns B..AC) of rows 130 and 131.
# ---------------------------------------------------------------------------
$row130 = @("6779694","Silkeborg IF","AGF Aarhus",0,1,"A",2.4,3.3,2.8,2.8,3.1,2.7,0,2,1.85,2.25,2.05,1.8,-1,-1,1.7,-1,0.8500000000000001,-1,0.8)
$row131 = @("6779696","Vejle","Odense BK",0,1,"A",2.5,3.2,2.8,2.7,3.2,2.7,0,1.95,1.9,2.25,1.925,1.925,-1,-1,1.7,-1,0.8999999999999999,-1,0.925)

for ($i = 0; $i -lt $cols.Count; $i++) {
    $ws.Cells.Item(130, $cols[$i]).Value = $row131[$i]
    $ws.Cells.Item(131, $cols[$i]).Value = $row130[$i]
}

# ---------------------------------------------------------------------------
# 3) Row 160 takes on the (refreshed) content of the former row 163 (match
#    7984010, Brondby vs Midtjylland) and row 161 takes on the refreshed
#    content of the former row 164 (match 7984011, FC Nordsjaelland vs AGF
#    Aarhus).
# ---------------------------------------------------------------------------
$ws.Cells.Item(160, 2).Value = 7984010
$ws.Cells.Item(160, 5).Value = 45403.54166666666
$ws.Cells.Item(160, 6).Value = "Brondby"
$ws.Cells.Item(160, 7).Value = "Midtjylland"
$ws.Cells.Item(160, 11).Value = 2.15
$ws.Cells.Item(160, 12).Value = 3.5
$ws.Cells.Item(160, 13).Value = 3.1
$ws.Cells.Item(160, 14).Value = 2.1
$ws.Cells.Item(160, 15).Value = 3.5
$ws.Cells.Item(160, 16).Value = 3.5
$ws.Cells.Item(160, 17).Value = -0.25
$ws.Cells.Item(160, 18).Value = 1.86
$ws.Cells.Item(160, 19).Value = 2.04
$ws.Cells.Item(160, 20).Value = 2.25
$ws.Cells.Item(160, 21).Value = 1.875
$ws.Cells.Item(160, 22).Value = 1.975

$ws.Cells.Item(161, 2).Value = 7984011
$ws.Cells.Item(161, 5).Value = 45404.58333333334
$ws.Cells.Item(161, 6).Value = "FC Nordsjaelland"
$ws.Cells.Item(161, 7).Value = "AGF Aarhus"
$ws.Cells.Item(161, 11).Value = 1.75
$ws.Cells.Item(161, 12).Value = 3.8
$ws.Cells.Item(161, 13).Value = 4
$ws.Cells.Item(161, 14).Value = 1.615
$ws.Cells.Item(161, 15).Value = 4
$ws.Cells.Item(161, 16).Value = 5.5
$ws.Cells.Item(161, 17).Value = -1
$ws.Cells.Item(161, 18).Value = 2.09
$ws.Cells.Item(161, 19).Value = 1.81
$ws.Cells.Item(161, 20).Value = 2.5
$ws.Cells.Item(161, 21).Value = 1.925
$ws.Cells.Item(161, 22).Value = 1.925

# ---------------------------------------------------------------------------
# 4) The matches that used to live in rows 160, 161, 162 are gone - delete
#    the (now obsolete) trailing rows 162:164.
# ---------------------------------------------------------------------------
$ws.Rows("162:164").Delete()
